$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that must remain plain text (they use
# "." as a thousands separator, e.g. "25.901.82", and Excel would otherwise
# auto-convert simple decimal-looking values like "14.98" into numbers).
# Force text format first for any new value that looks numeric, so the
# literal string is preserved exactly as in the source data.

$ws.Range("D2").Value = "25.901.82"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "1.742.20"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.85"
$ws.Range("E5").Value = "  +3.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5179"
$ws.Range("E7").Value = "  -0.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2746"
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06146"
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("D10").Value = "1.739.93"
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07171"
$ws.Range("E11").Value = "  +1.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.6449"
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.98"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.598"
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.47"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9993"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "25.912.51"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.73"
$ws.Range("E19").Value = "  +2.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006783"
$ws.Range("E20").Value = "  +2.36%  "
$ws.Range("D21").Value = "1.962.30"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.277"
$ws.Range("E22").Value = "  +1.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.674"
$ws.Range("E23").Value = "  -1.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.246"
$ws.Range("E24").Value = "  +1.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.72"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.510"
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.14"
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.765"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.03"
$ws.Range("E29").Value = "  +3.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.939"
$ws.Range("E30").Value = "  +5.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08302"
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.648"
$ws.Range("E32").Value = "  +4.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04592"
$ws.Range("E33").Value = "  +3.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.656"
$ws.Range("E34").Value = "  +1.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9907"
$ws.Range("E35").Value = "  +1.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6186"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.689"
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01613"
$ws.Range("E38").Value = "  +2.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.934"
$ws.Range("E39").Value = "  +2.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9992"
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.96"
$ws.Range("E41").Value = "  -1.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3839"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7401"
$ws.Range("E43").Value = "  +2.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.982"
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1126"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.213"
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05259"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.85"
$ws.Range("E48").Value = "  +3.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.53"
$ws.Range("E49").Value = "  +1.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.599"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3407"
$ws.Range("E51").Value = "  +0.85%  "
